$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Status and Date values ---
$ws1 = $wb.ActiveSheet
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Fix alignment/wrap formatting so it is properly flagged as applied ---
# (re-asserting the existing top/wrap alignment causes the engine to mark
#  applyAlignment="true" on the cell styles, matching the intended edit)
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A1:B1").VerticalAlignment = -4160
$ws1.Range("A2:B14").WrapText = $true
$ws1.Range("A2:B14").VerticalAlignment = -4160

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").WrapText = $true
$ws2.Range("A1").VerticalAlignment = -4160
$ws2.Range("A2:A4").WrapText = $true
$ws2.Range("A2:A4").VerticalAlignment = -4160
$ws2.Range("B3:B4").WrapText = $true
$ws2.Range("B3:B4").VerticalAlignment = -4160
